# Apply the updated cryptocurrency price/volume figures scraped on
# Mon Nov 13 02:48:30 UTC 2023 (see commit message). Only the cells whose
# text actually changed are touched; everything else is left as-is.
#
# "Price" (column D) and "Volume(1h)" (column E) are stored as plain text
# in this sheet, not numbers (e.g. "0.666", "  +1.24%  "). Some of the new
# Price strings look like ordinary numbers (no thousands separators), so a
# bare .Value assignment would let Excel silently reinterpret them as
# numeric values. IsText marks those cells so we force text formatting
# before writing, then restore the default style so no stray number
# format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Cells=@(@{Col='D'; Val='37.228.48'; IsText=0}, @{Col='E'; Val='  +0.75%  '; IsText=0})}
    @{Row=3; Cells=@(@{Col='D'; Val='2.063.66'; IsText=0}, @{Col='E'; Val='  +1.05%  '; IsText=0})}
    @{Row=4; Cells=@(@{Col='E'; Val='  -0.13%  '; IsText=0})}
    @{Row=5; Cells=@(@{Col='D'; Val='249.62'; IsText=1}, @{Col='E'; Val='  +0.83%  '; IsText=0})}
    @{Row=6; Cells=@(@{Col='D'; Val='0.666'; IsText=1}, @{Col='E'; Val='  +0.75%  '; IsText=0})}
    @{Row=7; Cells=@(@{Col='D'; Val='58.57'; IsText=1}, @{Col='E'; Val='  +5.80%  '; IsText=0})}
    @{Row=8; Cells=@(@{Col='E'; Val='  -0.07%  '; IsText=0})}
    @{Row=9; Cells=@(@{Col='D'; Val='0.386'; IsText=1}, @{Col='E'; Val='  +2.53%  '; IsText=0})}
    @{Row=10; Cells=@(@{Col='E'; Val='  +1.47%  '; IsText=0})}
    @{Row=11; Cells=@(@{Col='E'; Val='  +2.01%  '; IsText=0})}
    @{Row=12; Cells=@(@{Col='D'; Val='15.91'; IsText=1}, @{Col='E'; Val='  +1.01%  '; IsText=0})}
    @{Row=13; Cells=@(@{Col='D'; Val='0.921'; IsText=1}, @{Col='E'; Val='  +17.14%  '; IsText=0})}
    @{Row=14; Cells=@(@{Col='D'; Val='2.362.36'; IsText=0}, @{Col='E'; Val='  +0.97%  '; IsText=0})}
    @{Row=15; Cells=@(@{Col='E'; Val='  +4.55%  '; IsText=0})}
    @{Row=16; Cells=@(@{Col='D'; Val='2.068.01'; IsText=0}, @{Col='E'; Val='  +1.22%  '; IsText=0})}
    @{Row=17; Cells=@(@{Col='D'; Val='18.78'; IsText=1}, @{Col='E'; Val='  +14.51%  '; IsText=0})}
    @{Row=18; Cells=@(@{Col='D'; Val='37.247.04'; IsText=0}, @{Col='E'; Val='  +0.94%  '; IsText=0})}
    @{Row=19; Cells=@(@{Col='D'; Val='75.50'; IsText=1}, @{Col='E'; Val='  +2.63%  '; IsText=0})}
    @{Row=20; Cells=@(@{Col='D'; Val='0.0₃0915'; IsText=0}, @{Col='E'; Val='  +2.97%  '; IsText=0})}
    @{Row=21; Cells=@(@{Col='D'; Val='5.53'; IsText=1}, @{Col='E'; Val='  +4.44%  '; IsText=0})}
    @{Row=22; Cells=@(@{Col='D'; Val='239.46'; IsText=1}, @{Col='E'; Val='  +1.72%  '; IsText=0})}
    @{Row=23; Cells=@(@{Col='D'; Val='0.999'; IsText=1}, @{Col='E'; Val='  -0.19%  '; IsText=0})}
    @{Row=24; Cells=@(@{Col='E'; Val='  +6.05%  '; IsText=0})}
    @{Row=25; Cells=@(@{Col='D'; Val='2.22'; IsText=1}, @{Col='E'; Val='  +2.66%  '; IsText=0})}
    @{Row=26; Cells=@(@{Col='D'; Val='9.65'; IsText=1}, @{Col='E'; Val='  +6.69%  '; IsText=0})}
    @{Row=27; Cells=@(@{Col='D'; Val='171.88'; IsText=1}, @{Col='E'; Val='  +2.78%  '; IsText=0})}
    @{Row=28; Cells=@(@{Col='D'; Val='20.29'; IsText=1}, @{Col='E'; Val='  +3.30%  '; IsText=0})}
    @{Row=29; Cells=@(@{Col='E'; Val='  +19.28%  '; IsText=0})}
    @{Row=30; Cells=@(@{Col='E'; Val='  +1.31%  '; IsText=0})}
    @{Row=31; Cells=@(@{Col='E'; Val='  +6.03%  '; IsText=0})}
    @{Row=32; Cells=@(@{Col='D'; Val='4.88'; IsText=1}, @{Col='E'; Val='  +11.54%  '; IsText=0})}
    @{Row=33; Cells=@(@{Col='D'; Val='0.0631'; IsText=1}, @{Col='E'; Val='  +3.92%  '; IsText=0})}
    @{Row=34; Cells=@(@{Col='B'; Val='Kaspa'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; IsText=0}, @{Col='D'; Val='0.0887'; IsText=1}, @{Col='E'; Val='  +2.24%  '; IsText=0})}
    @{Row=35; Cells=@(@{Col='B'; Val='LidoDAOToken'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; IsText=0}, @{Col='D'; Val='2.35'; IsText=1}, @{Col='E'; Val='  +6.95%  '; IsText=0})}
    @{Row=36; Cells=@(@{Col='E'; Val='  +0.09%  '; IsText=0})}
    @{Row=37; Cells=@(@{Col='E'; Val='  +4.93%  '; IsText=0})}
    @{Row=38; Cells=@(@{Col='E'; Val='  +1.17%  '; IsText=0})}
    @{Row=39; Cells=@(@{Col='B'; Val='THORChain'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; IsText=0}, @{Col='D'; Val='5.20'; IsText=1}, @{Col='E'; Val='  +7.17%  '; IsText=0})}
    @{Row=40; Cells=@(@{Col='B'; Val='HuobiToken'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; IsText=0}, @{Col='D'; Val='3.13'; IsText=1}, @{Col='E'; Val='  -3.45%  '; IsText=0})}
    @{Row=41; Cells=@(@{Col='E'; Val='  -3.81%  '; IsText=0})}
    @{Row=42; Cells=@(@{Col='E'; Val='  +3.58%  '; IsText=0})}
    @{Row=43; Cells=@(@{Col='D'; Val='101.00'; IsText=1}, @{Col='E'; Val='  +6.37%  '; IsText=0})}
    @{Row=44; Cells=@(@{Col='B'; Val='InjectiveProtocol'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; IsText=0}, @{Col='D'; Val='17.64'; IsText=1}, @{Col='E'; Val='  +2.68%  '; IsText=0})}
    @{Row=45; Cells=@(@{Col='B'; Val='ARBITRUM'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; IsText=0}, @{Col='D'; Val='1.17'; IsText=1}, @{Col='E'; Val='  +5.74%  '; IsText=0})}
    @{Row=46; Cells=@(@{Col='E'; Val='  +1.27%  '; IsText=0})}
    @{Row=47; Cells=@(@{Col='B'; Val='FTXToken'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; IsText=0}, @{Col='D'; Val='3.96'; IsText=1}, @{Col='E'; Val='  +23.93%  '; IsText=0})}
    @{Row=48; Cells=@(@{Col='B'; Val='Maker'; IsText=0}, @{Col='C'; Val='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; IsText=0}, @{Col='D'; Val='1.311.99'; IsText=0}, @{Col='E'; Val='  +3.21%  '; IsText=0})}
    @{Row=49; Cells=@(@{Col='D'; Val='7.00'; IsText=1}, @{Col='E'; Val='  +5.37%  '; IsText=0})}
    @{Row=50; Cells=@(@{Col='E'; Val='  +1.56%  '; IsText=0})}
    @{Row=51; Cells=@(@{Col='D'; Val='2.250.23'; IsText=0}, @{Col='E'; Val='  +1.11%  '; IsText=0})}
)

foreach ($update in $updates) {
    foreach ($cell in $update.Cells) {
        $range = $ws.Range("$($cell.Col)$($update.Row)")
        if ($cell.IsText -eq 1) {
            $range.NumberFormat = "@"
            $range.Value = $cell.Val
            $range.Style = "Normal"
        } else {
            $range.Value = $cell.Val
        }
    }
}
